# Apply edits to the "Combat" sound-cue block on rows 21-32 of the
# Scarlet_Wald_M sheet, and update the sheet's current selection.
#
# Summary of the edit (reverse engineered from the target XML diff):
#  - Row 21's timing values / detail text are rewritten.
#  - Rows 22-24 (Angriff), 25-27 (Einstecken Leicht) get new
#    Start/End timestamps, and pick up the "black font" cell style
#    (xfId referenced by style index 1) on the C/D (Art/Details)
#    columns, matching the rows below them.
#  - Row 25 additionally becomes an "Einstecken Leicht" row (its
#    "Art" label flips from Angriff) and gains the time-format style
#    (style index 2) on its Start cell, just like row 21.
#  - Row 28 is cleared out entirely (content removed, formatting of
#    C28/D28 retained).
#  - Rows 29-30 (Einstecken Schwer) get new Start/End timestamps.
#  - Rows 31-32 are cleared out entirely (their formatting remains).
#  - The active sheet selection moves to E45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21 --------------------------------------------------------
$ws.Range("A21").Value = "00:22.0"
$ws.Range("B21").Value = "00:22.6"
$ws.Range("E21").Value = "Angriff"

# --- Row 22 ----------------------------------------------------------
$ws.Range("A22").Value = "00:22.7"
$ws.Range("B22").Value = "00:23.4"

# --- Row 23 ----------------------------------------------------------
$ws.Range("A23").Value = "00:23.6"
$ws.Range("B23").Value = "00:24.0"

# --- Row 24 ----------------------------------------------------------
$ws.Range("A24").Value = "00:24.2"
$ws.Range("B24").Value = "00:24.7"

# --- Row 25 (Angriff -> Einstecken Leicht) ---------------------------
$ws.Range("A25").Value = "00:28.8"
$ws.Range("B25").Value = "00:29.6"
$ws.Range("E25").Value = "Einstecken Leicht"

# --- Row 26 ------------------------------------------------------------
$ws.Range("A26").Value = "00:29.9"
$ws.Range("B26").Value = "00:30.5"

# --- Row 27 ------------------------------------------------------------
$ws.Range("A27").Value = "00:30.8"
$ws.Range("B27").Value = "00:31.7"

# --- Row 28: clear entirely (formatting of C28/D28 stays behind) -------
$ws.Range("A28").ClearContents()
$ws.Range("B28").ClearContents()
$ws.Range("C28").ClearContents()
$ws.Range("D28").ClearContents()
$ws.Range("E28").ClearContents()

# --- Row 29 (Einstecken Leicht -> Einstecken Schwer) --------------------
$ws.Range("A29").Value = "00:37.3"
$ws.Range("B29").Value = "00:38.5"
$ws.Range("E29").Value = "Einstecken Schwer"

# --- Row 30 ---------------------------------------------------------------
$ws.Range("A30").Value = "00:40.5"
$ws.Range("B30").Value = "00:41.9"

# --- Row 31: clear entirely (formatting stays behind) ----------------------
$ws.Range("A31").ClearContents()
$ws.Range("B31").ClearContents()
$ws.Range("C31").ClearContents()
$ws.Range("D31").ClearContents()
$ws.Range("E31").ClearContents()

# --- Row 32: clear entirely (formatting stays behind) -----------------------
$ws.Range("A32").ClearContents()
$ws.Range("B32").ClearContents()
$ws.Range("C32").ClearContents()
$ws.Range("D32").ClearContents()
$ws.Range("E32").ClearContents()

# --- Formatting: copy the "black font" style already used by C24/D24 -------
# onto C21:D23 so the whole Angriff block (rows 21-24) is consistently
# formatted, matching rows 25-27 below it.
$ws.Range("C24:D24").Copy()
$ws.Range("C21:D21").PasteSpecial(-4122)
$ws.Range("C22:D22").PasteSpecial(-4122)
$ws.Range("C23:D23").PasteSpecial(-4122)

# --- Formatting: copy the time-number-format style already used by A6 ------
# onto A21 and A25 (the Start cells that begin a new combat phase).
$ws.Range("A6").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A25").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Move the sheet's active selection --------------------------------------
$ws.Range("E45").Select()
